$d = $word.ActiveDocument

# The document ends with a bibliography entry, followed by a blank spacer
# paragraph, a page-break paragraph, and finally a site-generator "footer"
# paragraph (copyright / contact notice). That trailing footer block
# (the blank paragraph, the page-break paragraph, and the copyright
# paragraph itself) needs to be removed, leaving the bibliography text
# followed directly by the document's original closing blank + page-break
# paragraphs.

# Locate the copyright/footer paragraph by its distinctive text.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*luizeleno@usp.br*") {
        $target = $p
    }
}

if ($target -ne $null) {
    $idx = $target.Index
    $firstToRemove = $d.Paragraphs($idx - 2)
    $deleteRange = $d.Range($firstToRemove.Range.Start, $target.Range.End)
    $deleteRange.Delete()
}
